$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "OK" markers in columns B, D, F for rows 3-7 (per diff)
$ws.Range("B3").Value = "OK"
$ws.Range("D3").Value = "OK"
$ws.Range("F3").Value = "OK"

$ws.Range("B4").Value = "OK"
$ws.Range("D4").Value = "OK"
$ws.Range("F4").Value = "OK"

$ws.Range("B5").Value = "OK"
$ws.Range("D5").Value = "OK"
$ws.Range("F5").Value = "OK"

$ws.Range("B6").Value = "OK"
$ws.Range("D6").Value = "OK"
$ws.Range("F6").Value = "OK"

$ws.Range("B7").Value = "OK"

# New note row
$ws.Range("A11").Value = "Add Datalog for restart ESP"

# Update the selected cell to match the new state
$ws.Range("F7").Select()
